$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2686832740213523
$ws1.Range("C2").Value = 0.06378132118451026
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1199143468950749
$ws1.Range("F2").Value = 0.2540834845735027
$ws1.Range("G2").Value = 0.6391571553994733
$ws1.Range("H2").Value = 0.7873194221508828
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 411
$ws1.Range("K2").Value = 123
$ws1.Range("L2").Value = 0

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2303370786516854
$ws2.Range("D2").Value = 0.3744292237442922

$ws2.Range("B3").Value = 0.06378132118451026
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1199143468950749

$ws2.Range("B4").Value = 0.2686832740213523
$ws2.Range("C4").Value = 0.2686832740213523
$ws2.Range("D4").Value = 0.2686832740213523
$ws2.Range("E4").Value = 0.2686832740213523

$ws2.Range("B5").Value = 0.5318906605922551
$ws2.Range("C5").Value = 0.6151685393258427
$ws2.Range("D5").Value = 0.2471717853196836

$ws2.Range("B6").Value = 0.9533556530127515
$ws2.Range("C6").Value = 0.2686832740213523
$ws2.Range("D6").Value = 0.3617487672464665

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 123
$ws3.Range("C2").Value = 411

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
